# Refresh the crypto price/volume table with the latest scraped values.
# Source data stores every cell (coin name, link, price, 1h volume) as
# literal text, so numeric-looking prices must be written back as text
# too (not auto-coerced into Excel numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $cell = $ws.Range($ref)
    # Looks-like-a-number guard: force text formatting first so Excel
    # does not silently convert strings such as "224.37" into numbers,
    # then drop back to the default style so no formatting residue is
    # left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '34.062.89'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.811.75'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.09%  '
Set-TextCell "D5" '224.37'
$ws.Range("E5").Value = '  +0.13%  '
Set-TextCell "D6" '0.553'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.07%  '
Set-TextCell "D8" '31.72'
$ws.Range("E8").Value = '  -1.25%  '
Set-TextCell "D9" '0.288'
$ws.Range("E9").Value = '  +2.82%  '
Set-TextCell "D10" '0.0740'
$ws.Range("E10").Value = '  +12.80%  '
Set-TextCell "D11" '0.0930'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").Value = '2.075.32'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").Value = '1.819.78'
$ws.Range("E13").Value = '  +2.01%  '
Set-TextCell "D14" '10.80'
$ws.Range("E14").Value = '  -3.43%  '
Set-TextCell "D15" '0.639'
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '34.062.45'
$ws.Range("E16").Value = '  -0.47%  '
Set-TextCell "D17" '4.28'
$ws.Range("E17").Value = '  +1.89%  '
Set-TextCell "D18" '69.12'
$ws.Range("E18").Value = '  +0.55%  '
Set-TextCell "D19" '248.40'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = '0.0₃0792'
$ws.Range("E20").Value = '  +7.07%  '
Set-TextCell "D21" '10.98'
$ws.Range("E21").Value = '  +6.26%  '
Set-TextCell "D22" '1.00'
$ws.Range("E22").Value = '  +0.14%  '
Set-TextCell "D23" '4.22'
$ws.Range("E23").Value = '  +0.87%  '
Set-TextCell "D24" '2.15'
$ws.Range("E24").Value = '  +0.44%  '
Set-TextCell "D25" '159.79'
$ws.Range("E25").Value = '  +1.74%  '
Set-TextCell "D26" '16.51'
$ws.Range("E26").Value = '  +0.91%  '
Set-TextCell "D27" '7.17'
$ws.Range("E27").Value = '  +2.64%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  -0.01%  '
Set-TextCell "D30" '0.0527'
$ws.Range("E30").Value = '  +2.79%  '
Set-TextCell "D31" '3.74'
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("E32").Value = '  +1.81%  '
Set-TextCell "D33" '3.55'
$ws.Range("E33").Value = '  -0.75%  '
Set-TextCell "D34" '1.87'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").Value = '1.427.94'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("E36").Value = '  +0.47%  '
Set-TextCell "D37" '0.635'
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("E38").Value = '  +0.95%  '
Set-TextCell "D39" '0.953'
$ws.Range("E39").Value = '  +7.32%  '
$ws.Range("E40").Value = '  -1.67%  '
Set-TextCell "D41" '80.85'
$ws.Range("E41").Value = '  -2.35%  '
$ws.Range("E42").Value = '  +0.14%  '
Set-TextCell "D43" '2.14'
$ws.Range("E43").Value = '  +4.08%  '
Set-TextCell "D44" '6.03'
$ws.Range("E44").Value = '  +3.51%  '
Set-TextCell "D45" '1.05'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell "D46" '0.0495'
$ws.Range("E46").Value = '  -2.59%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.972.44'
$ws.Range("E47").Value = '  +1.53%  '
Set-TextCell "D48" '106.18'
$ws.Range("E48").Value = '  +8.02%  '
$ws.Range("E49").Value = '  -0.10%  '
Set-TextCell "D50" '11.75'
$ws.Range("E50").Value = '  -3.61%  '
$ws.Range("E51").Value = '  +4.02%  '
